$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(1, 15).Value = "combined lca db fungi"
$ws2.Cells.Item(2, 8).Value = "ALQSDSALK"
$ws2.Cells.Item(2, 9).Value = "AGPFGQLFR"
$ws2.Cells.Item(2, 10).Value = "AGPFGQLFR"
$ws2.Cells.Item(2, 15).Value = "Phialophora americana"
$ws2.Cells.Item(3, 8).Value = "AGPFGQLFR"
$ws2.Cells.Item(3, 9).Value = "GLSVGDGR"
$ws2.Cells.Item(3, 10).Value = "PGQLNSDLR"
$ws2.Cells.Item(3, 15).Value = "Pezizomycotina"
$ws2.Cells.Item(4, 8).Value = "AGPFGQLFRPD"
$ws2.Cells.Item(4, 9).Value = "VSHQVPR"
$ws2.Cells.Item(4, 15).Value = "Pezizomycotina"
$ws2.Cells.Item(5, 8).Value = "TLLDSVVEGK"
$ws2.Cells.Item(5, 9).Value = "TSGWFSK"
$ws2.Cells.Item(5, 15).Value = "Trichoderma gamsii"
$ws2.Cells.Item(6, 8).Value = "PGQLNSDLR"
$ws2.Cells.Item(6, 15).Value = "Scedosporium boydii"
$ws2.Cells.Item(7, 8).Value = "AAALQFTR"
$ws2.Cells.Item(7, 15).Value = "Endocarpon pusillum"
$ws2.Cells.Item(8, 8).Value = "AGPFGQLFRP"
$ws2.Cells.Item(8, 15).Value = "Fungi"
$ws2.Cells.Item(9, 8).Value = "VVTLLVNK"
$ws2.Cells.Item(9, 15).Value = "Rhodotorula"
$ws2.Cells.Item(10, 8).Value = "QLLLGFSK"
$ws2.Cells.Item(10, 15).Value = "Puccinia sorghi"
$ws2.Cells.Item(11, 8).Value = "SQEATLEK"
$ws2.Cells.Item(11, 15).Value = "Colletotrichum higginsianum"
$ws2.Cells.Item(12, 8).Value = "VSDTVVEPYNA"
$ws2.Cells.Item(12, 15).Value = "Dikarya"
$ws2.Cells.Item(13, 8).Value = "FYTTELDK"
$ws2.Cells.Item(13, 15).Value = "Mucor circinelloides f. circinelloides"
$ws2.Cells.Item(14, 8).Value = "LAFEPSNLK"
$ws2.Cells.Item(14, 15).Value = "Mortierella elongata"
$ws2.Cells.Item(15, 8).Value = "DSELCLR"
$ws2.Cells.Item(15, 15).Value = "Pezizomycotina"
$ws2.Cells.Item(16, 8).Value = "LTNTGSVK"
$ws2.Cells.Item(16, 15).Value = "Colletotrichum"
$ws2.Cells.Item(17, 8).Value = "VSDTVVEPYNATLSVHQLVEN"
$ws2.Cells.Item(17, 15).Value = "Pezizomycotina"
$ws2.Cells.Item(18, 8).Value = "AGPFGQLF"
$ws2.Cells.Item(18, 15).Value = "Pezizomycotina"
$ws2.Cells.Item(19, 8).Value = "QASLPLDR"
$ws2.Cells.Item(19, 15).Value = "Fonsecaea multimorphosa"
$ws2.Cells.Item(20, 15).Value = "Pezizomycotina"
$ws2.Cells.Item(21, 15).Value = "Laccaria amethystina"
$ws2.Cells.Item(22, 15).Value = "Suillus luteus"
$ws2.Cells.Item(23, 15).Value = "Stachybotrys"
$ws2.Cells.Item(24, 15).Value = "Pezizomycotina"
$ws2.Cells.Item(25, 15).Value = "Scedosporium boydii"

# Apply black-font style (matches existing style index used elsewhere in the workbook)
$ws2.Range("O2:O19").Font.Color = 0

# Make sheet2 the active/visible tab (moves tabSelected from sheet1 to sheet2
# and sets workbook bookViews activeTab)
$ws2.Activate()

# Final selection as last left by the author
$ws2.Range("O26").Select()
